# ------------------------------------------------------------------
# Gantt chart worksheet edit: insert a new "Toivoteet" column before
# the existing table, and append a new "Kommentti" column after it.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at A; existing columns A:E shift to B:F.
$ws.Columns.Item(1).Insert()

# --- Column A ("Toivoteet") + two corrected dates + one renamed task ---
$ws.Range("A1").Value = "Toivoteet"
$ws.Range("A2").Value = "Tulospalvelu"
$ws.Range("A3").Value = "Tulospalvelu"
$ws.Range("A5").Value = "Tulospalvelu"
$ws.Range("A6").Value = "Tulospalvelu"
$ws.Range("A7").Value = "Pelipaikka"
$ws.Range("A9").Value = "Pelipaikka"
$ws.Range("A11").Value = "Pelipaikka"
$ws.Range("A13").Value = "Pelipaikka"
$ws.Range("A14").Value = "Nettisivut"
$ws.Range("A15").Value = "Nettisivut"
$ws.Range("A16").Value = "Nettisivut"
$ws.Range("A18").Value = "Nettisivut"
$ws.Range("A19").Value = "Nettisivut"
$ws.Range("A20").Value = "Nettisivut"

# Two date corrections that are not pure column shifts
$ws.Range("D2").Value = 44197   # was 44105 (Loppunut for Tulospalvelun maarittely)
$ws.Range("C6").Value = 44378   # was 44317 (Alkaa for Aktiivinen kehitysvaihe)

# Row 20 task label changed to a new, distinct string
$ws.Range("B20").Value = "Aktiivinen kehitysvaihe II"

# --- Column G ("Kommentti") ---
$ws.Range("G1").Value = "Kommentti"
$ws.Range("G2").Value = "Ei kommenttia"
$ws.Range("G3").Value = "Ei kommenttia"
$ws.Range("G4").Value = "Tämä tehtävä epäonnistui"
$ws.Range("G5").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G6").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G7").Value = "Ei kommenttia"
$ws.Range("G8").Value = "Tämä tehtävä ei edennyt suunnitellusti, koska …"
$ws.Range("G9").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G10").Value = "Tämä tehtävä ei edennyt suunnitellusti, koska …"
$ws.Range("G11").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G12").Value = "Ei kommenttia"
$ws.Range("G13").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G14").Value = "Tämä tehtävä ei edennyt suunnitellusti, koska …"
$ws.Range("G15").Value = "Ei kommenttia"
$ws.Range("G16").Value = "Tämä tehtävä epäonnistui"
$ws.Range("G17").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G18").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G19").Value = "Tämä tehtävä ei ole vielä alkanut."
$ws.Range("G20").Value = "Tämä tehtävä ei edennyt suunnitellusti, koska …"

# Match the saved selection from the source workbook
$ws.Range("E2").Select()

Write-Output "Gantt column insert complete"
